# Split the single "join" check row into two rows: "join1" (visit key
# columns) and "join2" (crew key columns) so a template can describe
# multiple join rows instead of just one.

$wb = $excel.ActiveWorkbook

# --- "visit" sheet: the join row only needs to be renamed to "join1" ---
$wsVisit = $wb.Worksheets.Item("visit")
$wsVisit.Range("A7").Value = "join1"

# --- "crew" sheet: same rename ---
$wsCrew = $wb.Worksheets.Item("crew")
$wsCrew.Range("A6").Value = "join1"

# --- "count" sheet: rename existing join row to "join1" (drop the crew
#     marker from it) and add a new "join2" row carrying the crew marker ---
$wsCount = $wb.Worksheets.Item("count")
$wsCount.Range("A7").Value = "join1"
$wsCount.Range("F7").ClearContents()
$wsCount.Range("A8").Value = "join2"
$wsCount.Range("F8").Value = "crew"

# --- restore view state (zoom + selection) for each sheet, finishing on
#     the "count" sheet so it stays the active tab ---
$wsVisit.Activate()
$excel.ActiveWindow.Zoom = 150
$wsVisit.Range("A8").Select()

$wsCrew.Activate()
$excel.ActiveWindow.Zoom = 140
$wsCrew.Range("C13").Select()

$wsCount.Activate()
$wsCount.Range("C19").Select()
